$d = $word.ActiveDocument

# The sample "pom.xml.jim" reference in the body text was a typo / stale
# filename; it should read "pom.jam" (the jbang/sh template no longer
# needs the version, as explained in the commit message). Replace the
# text while keeping the run's existing character formatting
# (Courier New, 18pt, en-US) which Find/Replace preserves automatically.
$d.Content.Find.Execute("pom.xml.jim", $true, $false, $false, $false, $false,
                         $true, 1, $false, "pom.jam", 2)
